# Merge the split runs in the Title, Author and Abstract paragraphs into
# single runs each, without changing the visible text. Using Find/Execute
# over the containing range re-writes the matched range as one run, which
# collapses the previously separate "word" / "space" / "word" runs that made
# up each of these three paragraphs.

$d = $word.ActiveDocument

# Title: "Answers:" + " " + "Logarithms" -> "Answers: Logarithms"
$d.Content.Find.Execute("Answers: Logarithms", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Answers: Logarithms", 2)

# Author: "Zoë" + " " + "Gemmell" -> "Zoë Gemmell"
$d.Content.Find.Execute("Zoë Gemmell", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Zoë Gemmell", 2)

# Abstract: the many word/space runs -> one run with the full sentence
$d.Content.Find.Execute("Answers to questions relating to the study guide on logarithms.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Answers to questions relating to the study guide on logarithms.", 2)
